# Auto-generated edit script: updates market-price columns (H-N) on several
# Leve-profit sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4445195.5
$ws.Range("I137").Value = 702.75
$ws.Range("K137").Value = 2108.25
$ws.Range("M137").Value = 441.75
$ws.Range("H138").Value = 1348.6531
$ws.Range("I138").Value = 1161
$ws.Range("K138").Value = 3483
$ws.Range("M138").Value = 1657
$ws.Range("H140").Value = 70628.57000000001
$ws.Range("J140").Value = 70628.57000000001
$ws.Range("L140").Value = 70628.57000000001
$ws.Range("N140").Value = -80988.57000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 3004
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 3004
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 3004
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -3350
$ws.Range("H23").Value = 85000
$ws.Range("J23").Value = 85000
$ws.Range("L23").Value = 85000
$ws.Range("N23").Value = -85518
$ws.Range("H25").Value = 400
$ws.Range("I25").Value = 400
$ws.Range("K25").Value = 400
$ws.Range("M25").Value = 2
$ws.Range("H45").Value = 2429.4285
$ws.Range("I45").Value = 2430.2856
$ws.Range("J45").Value = 2428.5715
$ws.Range("K45").Value = 2430.2856
$ws.Range("L45").Value = 2428.5715
$ws.Range("M45").Value = -2053.2856
$ws.Range("N45").Value = -3182.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 944.7143
$ws.Range("I5").Value = 402
$ws.Range("J5").Value = 1668.3334
$ws.Range("K5").Value = 402
$ws.Range("L5").Value = 1668.3334
$ws.Range("M5").Value = -289
$ws.Range("N5").Value = -1894.3334
$ws.Range("H7").Value = 775.25
$ws.Range("I7").Value = 775.25
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 775.25
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -662.25
$ws.Range("N7").ClearContents()
$ws.Range("H11").Value = 2628.3333
$ws.Range("I11").Value = 1980
$ws.Range("J11").Value = 2952.5
$ws.Range("K11").Value = 1980
$ws.Range("L11").Value = 2952.5
$ws.Range("M11").Value = -1840
$ws.Range("N11").Value = -3232.5
$ws.Range("H12").Value = 438.75
$ws.Range("I12").Value = 438.75
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 438.75
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -270.75
$ws.Range("N12").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H24").Value = 2300
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 2300
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 2300
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -2770
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H134").Value = 2146.96
$ws.Range("I134").Value = 1400.6342
$ws.Range("J134").Value = 5546.8887
$ws.Range("K134").Value = 4201.902599999999
$ws.Range("L134").Value = 16640.6661
$ws.Range("M134").Value = -1666.902599999999
$ws.Range("N134").Value = -21710.6661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6414100.5
$ws.Range("I31").Value = 3797.6223
$ws.Range("J31").Value = 47623188
$ws.Range("K31").Value = 3797.6223
$ws.Range("L31").Value = 47623188
$ws.Range("M31").Value = -3502.6223
$ws.Range("N31").Value = -47623778
$ws.Range("H34").Value = 6414100.5
$ws.Range("I34").Value = 3797.6223
$ws.Range("J34").Value = 47623188
$ws.Range("K34").Value = 3797.6223
$ws.Range("L34").Value = 47623188
$ws.Range("M34").Value = -3595.6223
$ws.Range("N34").Value = -47623592

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 130.2
$ws.Range("I6").Value = 130.2
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 390.6
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -277.6
$ws.Range("N6").ClearContents()
$ws.Range("H7").Value = 150
$ws.Range("I7").Value = 158.33333
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 474.99999
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -362.99999
$ws.Range("N7").Value = -524
$ws.Range("H11").Value = 212
$ws.Range("I11").Value = 153.33333
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 459.99999
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = -319.99999
$ws.Range("N11").Value = -1180
$ws.Range("H12").Value = 96
$ws.Range("J12").Value = 94.40000000000001
$ws.Range("L12").Value = 283.2
$ws.Range("N12").Value = -629.2
$ws.Range("H13").Value = 316.625
$ws.Range("I13").Value = 212.75
$ws.Range("J13").Value = 420.5
$ws.Range("K13").Value = 638.25
$ws.Range("L13").Value = 1261.5
$ws.Range("M13").Value = -470.25
$ws.Range("N13").Value = -1597.5
$ws.Range("H23").Value = 94.55
$ws.Range("I23").Value = 32.75
$ws.Range("J23").Value = 110
$ws.Range("K23").Value = 98.25
$ws.Range("L23").Value = 330
$ws.Range("M23").Value = 136.75
$ws.Range("N23").Value = -800
$ws.Range("H25").Value = 1070.2
$ws.Range("I25").Value = 240
$ws.Range("J25").Value = 1900.4
$ws.Range("K25").Value = 720
$ws.Range("L25").Value = 5701.200000000001
$ws.Range("M25").Value = -551
$ws.Range("N25").Value = -6039.200000000001
$ws.Range("H30").Value = 1070.2
$ws.Range("I30").Value = 240
$ws.Range("J30").Value = 1900.4
$ws.Range("K30").Value = 720
$ws.Range("L30").Value = 5701.200000000001
$ws.Range("M30").Value = -618
$ws.Range("N30").Value = -5905.200000000001
$ws.Range("H39").Value = 540.4595
$ws.Range("J39").Value = 540.4595
$ws.Range("L39").Value = 1621.3785
$ws.Range("N39").Value = -2209.3785
$ws.Range("H46").Value = 561.1111
$ws.Range("I46").Value = 592.8570999999999
$ws.Range("J46").Value = 450
$ws.Range("K46").Value = 1778.5713
$ws.Range("L46").Value = 1350
$ws.Range("M46").Value = -1687.5713
$ws.Range("N46").Value = -1532
$ws.Range("H132").Value = 957.8
$ws.Range("I132").Value = 808.3889
$ws.Range("J132").Value = 1181.9166
$ws.Range("K132").Value = 7275.5001
$ws.Range("L132").Value = 10637.2494
$ws.Range("M132").Value = -4745.5001
$ws.Range("N132").Value = -15697.2494

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 902
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 902
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 902
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -1180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 500
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 500
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 500
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -960
$ws.Range("H46").Value = 726.95654
$ws.Range("I46").Value = 816.6667
$ws.Range("J46").Value = 695.2941
$ws.Range("K46").Value = 816.6667
$ws.Range("L46").Value = 695.2941
$ws.Range("M46").Value = -628.6667
$ws.Range("N46").Value = -1071.2941
$ws.Range("H61").Value = 1500.0588
$ws.Range("I61").Value = 1397.6428
$ws.Range("J61").Value = 1978
$ws.Range("K61").Value = 1397.6428
$ws.Range("L61").Value = 1978
$ws.Range("M61").Value = -1195.6428
$ws.Range("N61").Value = -2382
$ws.Range("H113").Value = 1500.0588
$ws.Range("I113").Value = 1397.6428
$ws.Range("J113").Value = 1978
$ws.Range("K113").Value = 1397.6428
$ws.Range("L113").Value = 1978
$ws.Range("M113").Value = 772.3571999999999
$ws.Range("N113").Value = -6318

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
